# Apply the schema.xlsx update ("Add files via upload"):
#   - sections : rename the CMKS group's "Discap" SectionKey to "Discap_CMKS"
#   - fields   : add the new "Discap_CMKS" field definition row (select/radio)
#   - options  : add the "Discap_CMKS_opts" option list (D1200.."Diğer")
#   - leave the workbook positioned/selected the way it was left after saving
#
# New-value cells are written in the same order Excel would first encounter
# them (options sheet row 173 first) so the two brand-new shared strings
# ("Discap_CMKS_opts", "Discap_CMKS") land at the expected table indices.

$wb = $excel.ActiveWorkbook

$wsSections = $wb.Worksheets.Item("sections")
$wsFields   = $wb.Worksheets.Item("fields")
$wsOptions  = $wb.Worksheets.Item("options")

# --- options: new "Discap_CMKS_opts" list, rows 173-178 ---
$wsOptions.Range("A173").Value = "Discap_CMKS_opts"
$wsOptions.Range("B173").Value = "D1200"
$wsOptions.Range("C173").Value = "(1200 mm rulo dış çapı)"
$wsOptions.Range("D173").Value = 1

$wsOptions.Range("A174").Value = "Discap_CMKS_opts"
$wsOptions.Range("B174").Value = "D1300"
$wsOptions.Range("C174").Value = "(1300 mm rulo dış çapı)"
$wsOptions.Range("D174").Value = 2

$wsOptions.Range("A175").Value = "Discap_CMKS_opts"
$wsOptions.Range("B175").Value = "D1600"
$wsOptions.Range("C175").Value = "(1600 mm rulo dış çapı)"
$wsOptions.Range("D175").Value = 3

$wsOptions.Range("A176").Value = "Discap_CMKS_opts"
$wsOptions.Range("B176").Value = "D1800"
$wsOptions.Range("C176").Value = "(1800 mm rulo dış çapı)"
$wsOptions.Range("D176").Value = 4

$wsOptions.Range("A177").Value = "Discap_CMKS_opts"
$wsOptions.Range("B177").Value = "D2000"
$wsOptions.Range("C177").Value = "(2000 mm rulo dış çapı)"
$wsOptions.Range("D177").Value = 5

$wsOptions.Range("A178").Value = "Discap_CMKS_opts"
$wsOptions.Range("B178").Value = "Diğer"
$wsOptions.Range("C178").Value = "Lütfen aşağıdaki alana değer giriniz"
$wsOptions.Range("D178").Value = 6

# --- sections: row 37, column D (SectionKey) Discap -> Discap_CMKS ---
$wsSections.Range("D37").Value = "Discap_CMKS"

# --- fields: new "Discap_CMKS" field definition row, row 60 ---
$wsFields.Range("A60").Value = "Discap_CMKS"
$wsFields.Range("B60").Value = "Discap_CMKS"
$wsFields.Range("C60").Value = "Rulo Dış Çapı"
$wsFields.Range("D60").Value = "select"
$wsFields.Range("E60").Value = $true
$wsFields.Range("F60").Value = "Discap_CMKS_opts"
$wsFields.Range("Q60").Value = "radio"

# --- final selection / active-sheet bookkeeping ---
$wsSections.Range("D37").Select()
$wsOptions.Range("A173:A178").Select()
$wsFields.Range("O60").Select()
$wsFields.Activate()
